$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1139.1666
$ws.Range("J17").Value = 1139.1666
$ws.Range("L17").Value = 3417.4998
$ws.Range("N17").Value = -3753.4998

$ws.Range("H42").Value = 202.16667
$ws.Range("I42").Value = 68.666664
$ws.Range("J42").Value = 335.66666
$ws.Range("K42").Value = 205.999992
$ws.Range("L42").Value = 1006.99998
$ws.Range("M42").Value = 24.00000800000001
$ws.Range("N42").Value = -1466.99998

$ws.Range("H138").Value = 2416.0862
$ws.Range("J138").Value = 3641.4666
$ws.Range("L138").Value = 10924.3998
$ws.Range("N138").Value = -21204.3998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4019.64
$ws.Range("I32").Value = 3907.796
$ws.Range("J32").Value = 9500
$ws.Range("K32").Value = 3907.796
$ws.Range("L32").Value = 9500
$ws.Range("M32").Value = -3620.796
$ws.Range("N32").Value = -10074

$ws.Range("H110").Value = 1630
$ws.Range("I110").Value = 1610.2
$ws.Range("J110").Value = 1663
$ws.Range("K110").Value = 1610.2
$ws.Range("L110").Value = 1663
$ws.Range("M110").Value = 434.8
$ws.Range("N110").Value = -5753

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 803.1667
$ws.Range("I16").Value = 772.5789
$ws.Range("J16").Value = 919.4
$ws.Range("K16").Value = 772.5789
$ws.Range("L16").Value = 919.4
$ws.Range("M16").Value = -485.5789
$ws.Range("N16").Value = -1493.4

$ws.Range("H31").Value = 2059.2415
$ws.Range("I31").Value = 1145.2
$ws.Range("J31").Value = 4090.4443
$ws.Range("K31").Value = 1145.2
$ws.Range("L31").Value = 4090.4443
$ws.Range("M31").Value = -850.2
$ws.Range("N31").Value = -4680.4443

$ws.Range("H34").Value = 2059.2415
$ws.Range("I34").Value = 1145.2
$ws.Range("J34").Value = 4090.4443
$ws.Range("K34").Value = 1145.2
$ws.Range("L34").Value = 4090.4443
$ws.Range("M34").Value = -943.2
$ws.Range("N34").Value = -4494.4443

$ws.Range("H68").Value = 18125.7
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 18125.7
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 18125.7
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -19623.7

$ws.Range("H71").Value = 18125.7
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 18125.7
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 54377.10000000001
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -61865.10000000001

$ws.Range("H107").Value = 1494.875
$ws.Range("I107").Value = 1967.0667
$ws.Range("J107").Value = 707.8889
$ws.Range("K107").Value = 1967.0667
$ws.Range("L107").Value = 707.8889
$ws.Range("M107").Value = -47.06670000000008
$ws.Range("N107").Value = -4547.8889

$ws.Range("H113").Value = 803.1667
$ws.Range("I113").Value = 772.5789
$ws.Range("J113").Value = 919.4
$ws.Range("K113").Value = 772.5789
$ws.Range("L113").Value = 919.4
$ws.Range("M113").Value = 1397.4211
$ws.Range("N113").Value = -5259.4

$ws.Range("H122").Value = 3476240
$ws.Range("I122").Value = 4634415
$ws.Range("J122").Value = 1715.5555
$ws.Range("K122").Value = 13903245
$ws.Range("L122").Value = 5146.666499999999
$ws.Range("M122").Value = -13900795
$ws.Range("N122").Value = -10046.6665

$ws.Range("H132").Value = 1550.6511
$ws.Range("J132").Value = 2848
$ws.Range("L132").Value = 8544
$ws.Range("N132").Value = -13604

$ws.Range("H134").Value = 2294.5898
$ws.Range("I134").Value = 2006.5625
$ws.Range("J134").Value = 3611.2856
$ws.Range("K134").Value = 6019.6875
$ws.Range("L134").Value = 10833.8568
$ws.Range("M134").Value = -3484.6875
$ws.Range("N134").Value = -15903.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 758.3953
$ws.Range("I5").Value = 592.48
$ws.Range("J5").Value = 988.8333
$ws.Range("K5").Value = 1777.44
$ws.Range("L5").Value = 2966.4999
$ws.Range("M5").Value = -1665.44
$ws.Range("N5").Value = -3190.4999

$ws.Range("H103").Value = 4250518
$ws.Range("I103").Value = 4250518
$ws.Range("K103").Value = 12751554
$ws.Range("M103").Value = -12750675

$ws.Range("H107").Value = 100385.85
$ws.Range("I107").Value = 77283
$ws.Range("J107").Value = 143291.14
$ws.Range("K107").Value = 231849
$ws.Range("L107").Value = 429873.42
$ws.Range("M107").Value = -229929
$ws.Range("N107").Value = -433713.42

$ws.Range("H113").Value = 610.36365
$ws.Range("I113").Value = 608.0417
$ws.Range("J113").Value = 616.55554
$ws.Range("K113").Value = 1824.1251
$ws.Range("L113").Value = 1849.66662
$ws.Range("M113").Value = 345.8749
$ws.Range("N113").Value = -6189.66662

$ws.Range("H131").Value = 3041.9185
$ws.Range("J131").Value = 1927.6957
$ws.Range("L131").Value = 5783.0871
$ws.Range("N131").Value = -15863.0871

$ws.Range("H132").Value = 40000996
$ws.Range("I132").Value = 66667480
$ws.Range("J132").Value = 1269
$ws.Range("K132").Value = 600007320
$ws.Range("L132").Value = 11421
$ws.Range("M132").Value = -600004790
$ws.Range("N132").Value = -16481

$ws.Range("H135").Value = 758.3953
$ws.Range("I135").Value = 592.48
$ws.Range("J135").Value = 988.8333
$ws.Range("K135").Value = 5332.32
$ws.Range("L135").Value = 8899.4997
$ws.Range("M135").Value = -2797.32
$ws.Range("N135").Value = -13969.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2788.889
$ws.Range("I80").Value = 2933.3333
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 2933.3333
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -1935.3333
$ws.Range("N80").Value = -4496

$ws.Range("H83").Value = 2788.889
$ws.Range("I83").Value = 2933.3333
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 14666.6665
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -9674.666499999999
$ws.Range("N83").Value = -22484

$ws.Range("H102").Value = 2831.1428
$ws.Range("I102").Value = 2701.5
$ws.Range("J102").Value = 3004
$ws.Range("K102").Value = 2701.5
$ws.Range("L102").Value = 3004
$ws.Range("M102").Value = -1079.5
$ws.Range("N102").Value = -6248

$ws.Range("H132").Value = 2288.7646
$ws.Range("I132").Value = 1478.6666
$ws.Range("J132").Value = 3200.125
$ws.Range("K132").Value = 4435.9998
$ws.Range("L132").Value = 9600.375
$ws.Range("M132").Value = -1905.9998
$ws.Range("N132").Value = -14660.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 368.27777
$ws.Range("I22").Value = 375.66666
$ws.Range("J22").Value = 360.8889
$ws.Range("K22").Value = 375.66666
$ws.Range("L22").Value = 360.8889
$ws.Range("M22").Value = -80.66665999999998
$ws.Range("N22").Value = -950.8888999999999

$ws.Range("H27").Value = 368.27777
$ws.Range("I27").Value = 375.66666
$ws.Range("J27").Value = 360.8889
$ws.Range("K27").Value = 375.66666
$ws.Range("L27").Value = 360.8889
$ws.Range("M27").Value = -268.66666
$ws.Range("N27").Value = -574.8888999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3846.724
$ws.Range("I132").Value = 4521.3257
$ws.Range("J132").Value = 1912.8667
$ws.Range("K132").Value = 13563.9771
$ws.Range("L132").Value = 5738.6001
$ws.Range("M132").Value = -11033.9771
$ws.Range("N132").Value = -10798.6001
